$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in A2 (new shared string "Alteração no xlsx2")
$ws.Range("A2").Value = "Alteração no xlsx2"

# Autofit column A to match the bestFit width seen in the diff
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Select A3 as the active cell, like in the final sheetView selection
$ws.Range("A3").Select() | Out-Null
